$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

# Row 2
$ws.Range("D2").Value = '51.931.39'
$ws.Range("E2").Value = '  +0.65%  '

# Row 3
$ws.Range("D3").Value = '2.788.36'
$ws.Range("E3").Value = '  -0.73%  '

# Row 4
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
Set-TextValue $ws.Range("D5") '358.42'
$ws.Range("E5").Value = '  +1.84%  '

# Row 6
Set-TextValue $ws.Range("D6") '109.70'
$ws.Range("E6").Value = '  -2.73%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.565'
$ws.Range("E7").Value = '  +0.58%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.999'
$ws.Range("E8").Value = '  +0.08%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.593'
$ws.Range("E9").Value = '  -0.24%  '

# Row 10
Set-TextValue $ws.Range("D10") '40.03'
$ws.Range("E10").Value = '  -2.90%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.0853'
$ws.Range("E11").Value = '  +0.35%  '

# Row 12
Set-TextValue $ws.Range("D12") '0.133'
$ws.Range("E12").Value = '  +1.40%  '

# Row 13
Set-TextValue $ws.Range("D13") '19.51'
$ws.Range("E13").Value = '  -1.66%  '

# Row 14
Set-TextValue $ws.Range("D14") '7.60'
$ws.Range("E14").Value = '  -1.23%  '

# Row 15
$ws.Range("D15").Value = '3.228.03'
$ws.Range("E15").Value = '  -0.59%  '

# Row 16
$ws.Range("D16").Value = '2.784.32'
$ws.Range("E16").Value = '  -0.65%  '

# Row 17
Set-TextValue $ws.Range("D17") '0.946'

# Row 18
$ws.Range("D18").Value = '51.897.46'
$ws.Range("E18").Value = '  +1.19%  '

# Row 19
Set-TextValue $ws.Range("D19") '7.41'
$ws.Range("E19").Value = '  +0.02%  '

# Row 20
$ws.Range("E20").Value = '  -1.53%  '

# Row 21
Set-TextValue $ws.Range("D21") '13.04'
$ws.Range("E21").Value = '  -1.93%  '

# Row 22
$ws.Range("D22").Value = '0.0₃0982'
$ws.Range("E22").Value = '  -0.70%  '

# Row 23
Set-TextValue $ws.Range("D23") '274.21'
$ws.Range("E23").Value = '  +1.50%  '

# Row 24
Set-TextValue $ws.Range("D24") '70.23'
$ws.Range("E24").Value = '  +1.14%  '

# Row 25
Set-TextValue $ws.Range("D25") '2.74'
$ws.Range("E25").Value = '  -0.26%  '

# Row 26
Set-TextValue $ws.Range("D26") '26.69'
$ws.Range("E26").Value = '  +0.36%  '

# Row 27
$ws.Range("E27").Value = '  -0.10%  '

# Row 28
Set-TextValue $ws.Range("D28") '10.19'
$ws.Range("E28").Value = '  -0.56%  '

# Row 29
$ws.Range("E29").Value = '  +4.80%  '

# Row 30
$ws.Range("E30").Value = '  -0.92%  '

# Row 31
Set-TextValue $ws.Range("D31") '0.0465'
$ws.Range("E31").Value = '  +4.65%  '

# Row 32
Set-TextValue $ws.Range("D32") '51.55'
$ws.Range("E32").Value = '  +1.98%  '

# Row 33
Set-TextValue $ws.Range("D33") '34.41'
$ws.Range("E33").Value = '  +1.71%  '

# Row 34
Set-TextValue $ws.Range("D34") '5.73'
$ws.Range("E34").Value = '  -1.76%  '

# Row 35
Set-TextValue $ws.Range("D35") '0.0844'
$ws.Range("E35").Value = '  +3.11%  '

# Row 36
Set-TextValue $ws.Range("D36") '5.28'
$ws.Range("E36").Value = '  +4.97%  '

# Row 37
$ws.Range("E37").Value = '  +0.26%  '

# Row 38
Set-TextValue $ws.Range("D38") '3.23'
$ws.Range("E38").Value = '  +1.28%  '

# Row 39
$ws.Range("E39").Value = '  -2.15%  '

# Row 40
Set-TextValue $ws.Range("D40") '17.99'
$ws.Range("E40").Value = '  +0.10%  '

# Row 41
Set-TextValue $ws.Range("D41") '2.55'
$ws.Range("E41").Value = '  +2.08%  '

# Row 42
$ws.Range("E42").Value = '  -1.10%  '

# Row 45
Set-TextValue $ws.Range("D45") '22.06'
$ws.Range("E45").Value = '  -6.50%  '

# Row 46
$ws.Range("D46").Value = '2.073.79'
$ws.Range("E46").Value = '  +0.19%  '

# Row 47
Set-TextValue $ws.Range("D47") '3.25'
$ws.Range("E47").Value = '  -1.43%  '

# Row 48
$ws.Range("E48").Value = '  -4.96%  '

# Row 49
Set-TextValue $ws.Range("D49") '5.74'
$ws.Range("E49").Value = '  +2.06%  '

# Row 50
Set-TextValue $ws.Range("D50") '0.931'
$ws.Range("E50").Value = '  +1.53%  '

# Row 51
Set-TextValue $ws.Range("D51") '8.95'
$ws.Range("E51").Value = '  +0.78%  '

# Row 43 and 44 swap (Monero <-> WEMIXToken)
$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D43") '2.25'
$ws.Range("E43").Value = '  -1.59%  '

$ws.Range("B44").Value = 'Monero'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D44") '121.88'
$ws.Range("E44").Value = '  -3.49%  '
